$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Add new row 13 with translations for "Challenge" (mirrors header columns
# A=English, B=Francais, C=Arabic, D=Detche, E=Spanish, F=italian)
# Write B before A so the shared-string table records DEFI before CHALLENGE,
# matching the order new strings were appended upstream.
$ws.Range("B13").Value = "DEFI"
$ws.Range("A13").Value = "CHALLENGE"
$ws.Range("C13").Value = "تحد"
$ws.Range("D13").Value = "UITDAGING"
$ws.Range("E13").Value = "DESAFÍO"
$ws.Range("F13").Value = "SFIDA"

# Update selection/view to match the post-edit state (scrolled down one row,
# active cell moved to the newly added F13)
$ws.Range("F13").Select()
$excel.ActiveWindow.ScrollRow = 2
$excel.ActiveWindow.ScrollColumn = 1
